# Asthma_IR.xlsx edit script
# - Excludes Alaska (FIPS 2) and Hawaii (FIPS 15) from the yearly incidence/
#   prevalence tables (clears their computed columns) in the 2006-2010 sheets.
# - Removes the trailing placeholder rows (FIPS 66/72/78 etc.) that were left
#   over at the bottom of each yearly sheet.
# - Rebuilds the Aggregate sheet: drops Alaska and Hawaii entirely and
#   refreshes the recalculated incidence numbers for the remaining states.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Yearly sheets: 2006
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2006")
$ws.Range("C3:E3").ClearContents()
$ws.Range("C13:E13").ClearContents()
$ws.Rows("53:54").Delete()

# ---------------------------------------------------------------------------
# Yearly sheets: 2007
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2007")
$ws.Range("C3:E3").ClearContents()
$ws.Range("C13:E13").ClearContents()
$ws.Rows("53:55").Delete()

# ---------------------------------------------------------------------------
# Yearly sheets: 2008
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2008")
$ws.Range("C13:E13").ClearContents()
$ws.Rows("53:55").Delete()

# ---------------------------------------------------------------------------
# Yearly sheets: 2009
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2009")
$ws.Range("C13:E13").ClearContents()
$ws.Rows("53:55").Delete()

# ---------------------------------------------------------------------------
# Yearly sheets: 2010
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2010")
$ws.Range("C13:E13").ClearContents()
$ws.Rows("53:55").Delete()

# ---------------------------------------------------------------------------
# Aggregate sheet: drop Alaska (row 2) and Hawaii (row 7 once Alaska is gone),
# then push the refreshed incidence figures into the remaining rows.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Aggregate")
$ws.Rows("2:2").Delete()
$ws.Rows("7:7").Delete()

$updates = @(
    @{ Row = 2;  C = 42621.918836024626;  D = 2802421.7470672023;  E = 15.208959493919652 },
    @{ Row = 3;  C = 156599.0264181489;   D = 16850453.394726697;  E = 9.29346070101331 },
    @{ Row = 5;  C = 3184.446344618885;   D = 179492.85599165817;  E = 17.74135425627681 },
    @{ Row = 6;  C = 94786.03809289802;   D = 10458073.543508092;  E = 9.063431969431596 },
    @{ Row = 8;  C = 105218.93189138013;  D = 6936762.371679282;   E = 15.168305652354107 },
    @{ Row = 9;  C = 11510.486502884885;  D = 1829733.9661702944;  E = 6.290797851327419 },
    @{ Row = 10; C = 27509.487892604753;  D = 3059760.3617240055;  E = 8.990732815789757 },
    @{ Row = 12; C = 6662.149701015229;   D = 722763.1937653258;   E = 9.217610634415294 },
    @{ Row = 13; C = 64870.5736086545;    D = 5816583.904105316;   E = 11.152692831073782 },
    @{ Row = 14; C = 126101.97053771566;  D = 10491065.363787048;  E = 12.0199394594369 },
    @{ Row = 16; C = 46410.034572168624;  D = 3600272.063188931;   E = 12.890702079625912 },
    @{ Row = 17; C = 3295.889765187216;   D = 768012.0066596719;   E = 4.291456040540417 },
    @{ Row = 18; C = 18262.226573581298;  D = 2014605.1364890952;  E = 9.064916118206348 },
    @{ Row = 19; C = 9423.24961501254;    D = 788301.5634397555;   E = 11.953863917120971 },
    @{ Row = 22; C = 221226.36562273267;  D = 15027480.545143578;  E = 14.721454135851552 },
    @{ Row = 28; C = 381999.04549088026;  D = 22992022.96837785;   E = 16.61441648767765 },
    @{ Row = 30; C = 6498.410650110485;   D = 563279.5559912089;   E = 11.536741536225586 },
    @{ Row = 33; C = 14403.531708608236;  D = 1174446.719861135;   E = 12.264099737373606 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
}
